$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E7").Value = 11.96139999999999
$ws.Range("A10").Value = -20.45559999999997
$ws.Range("A12").Value = -22.54360000000004
$ws.Range("B13").Value = 6.346199999999995
$ws.Range("A18").Value = -22.45490000000003
$ws.Range("E20").Value = 13.2238
